$wb = $excel.ActiveWorkbook

$wsJob = $wb.Worksheets.Item("Job to Run")
$wsAll = $wb.Worksheets.Item("All")

# --- "Job to Run" sheet: point the template job at the new batch-params file ---
# (this mints the new filename shared string first -> index 11)
$wsJob.Range("E2").Value = "MZmine3_batch_params_LCMSMS_HE_for_Commandline_2024_9_auto_GNPS_run.xml"

# --- "All" sheet: append a new GNPS auto-run template block (rows 11-12) ---
$wsAll.Range("A11").Value = "Job Name"
$wsAll.Range("B11").Value = "EXP num replicates"
$wsAll.Range("C11").Value = "CTRL num replicates"
$wsAll.Range("D11").Value = "Ionization"
$wsAll.Range("E11").Value = "MZmine3 batch template"

$wsAll.Range("A12").Value = "Anid_HE_TJGIp11_pos_2018"
$wsAll.Range("B12").Value = 3
$wsAll.Range("C12").Value = 3
$wsAll.Range("D12").Value = "POS"
$wsAll.Range("E12").Value = "MZmine3_batch_params_LCMSMS_HE_for_Commandline_2024_9_auto_GNPS_run.xml"

$wsAll.Range("A12:D12").VerticalAlignment = -4108

# Highlighted header label above the new block
# (mints the "Auto-run GNPS Job" shared string last -> index 12)
$wsAll.Range("A10").Value = "Auto-run GNPS Job"
$wsAll.Range("A10").Interior.Color = 65535

# --- Selections / active sheet: "All" becomes the active tab ---
$wsJob.Activate()
$wsJob.Range("A2:E2").Select()

$wsAll.Activate()
$wsAll.Range("A7").Select()
